$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 111; this shifts all existing rows
# 111-196 down to 112-197 and updates the used range to A1:T197.
$ws.Rows.Item(111).Insert()

# Populate the newly inserted row 111 with the new data record.
$ws.Cells.Item(111, 1).Value = 4
$ws.Cells.Item(111, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(111, 3).Value = "Los Lagos"
$ws.Cells.Item(111, 4).Value = 44673
$ws.Cells.Item(111, 5).Value = 10
$ws.Cells.Item(111, 6).Value = "Fruta"
$ws.Cells.Item(111, 7).Value = 100102
$ws.Cells.Item(111, 8).Value = "Cítricos"
$ws.Cells.Item(111, 9).Value = 100102004
$ws.Cells.Item(111, 10).Value = "Mandarina"
$ws.Cells.Item(111, 11).Value = "Murcott"
$ws.Cells.Item(111, 12).Value = "Primera"
$ws.Cells.Item(111, 13).Value = 700
$ws.Cells.Item(111, 14).Value = 13000
$ws.Cells.Item(111, 15).Value = 14000
$ws.Cells.Item(111, 16).Value = 13500
$ws.Cells.Item(111, 17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(111, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(111, 19).Value = 1350
$ws.Cells.Item(111, 20).Value = 10
